$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) - labels unchanged, but rewritten for completeness ---
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# --- Ensure A2:A7 are treated as plain text (avoid auto date conversion) ---
$ws.Range("A2:A7").NumberFormat = "@"

# --- Row labels (A2:A7) ---
$ws.Range("A2").Value = "2025-09-30"
$ws.Range("A3").Value = "2025-10-15"
$ws.Range("A4").Value = "2025-10-30"
$ws.Range("A5").Value = "2025-11-15"
$ws.Range("A6").Value = "2025-11-30"
$ws.Range("A7").Value = "2025-12-15"

# --- Restore default (Normal) style so formatting matches the rest of the sheet ---
$ws.Range("A2:A7").Style = "Normal"

# --- Data values (B2:K7) ---
$ws.Range("B2").Value = 0.2113612845817191
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 0.0
$ws.Range("E2").Value = 0.0
$ws.Range("F2").Value = 0.0
$ws.Range("G2").Value = 0.0
$ws.Range("H2").Value = 0.0
$ws.Range("I2").Value = 0.0
$ws.Range("J2").Value = 0.0
$ws.Range("K2").Value = 0.0
$ws.Range("B3").Value = 0.22846811570041856
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 0.0022632260226792948
$ws.Range("E3").Value = 0.0013952951316078875
$ws.Range("F3").Value = 0.0019083025805341702
$ws.Range("G3").Value = 0.0010010192655974803
$ws.Range("H3").Value = -0.00003627697265169035
$ws.Range("I3").Value = -0.000495687437614611
$ws.Range("J3").Value = 0.0
$ws.Range("K3").Value = 0.005085330807441113
$ws.Range("B4").Value = 0.34127803299941145
$ws.Range("C4").Value = 0.05591714355495465
$ws.Range("D4").Value = 0.0
$ws.Range("E4").Value = -0.0008143886754195167
$ws.Range("F4").Value = -0.000009371645147495026
$ws.Range("G4").Value = 0.0
$ws.Range("H4").Value = 0.000582617554790624
$ws.Range("I4").Value = -0.0023759055295641505
$ws.Range("J4").Value = 0.0029976442716891814
$ws.Range("K4").Value = -0.00011946344348232052
$ws.Range("B5").Value = 0.3597433531423276
$ws.Range("C5").Value = 0.0
$ws.Range("D5").Value = 0.00014699654210745616
$ws.Range("E5").Value = -0.0003673714797570992
$ws.Range("F5").Value = 0.010878750558187421
$ws.Range("G5").Value = -0.0024531109789531075
$ws.Range("H5").Value = 0.00034354718127220215
$ws.Range("I5").Value = -0.000621280498114726
$ws.Range("J5").Value = 0.0
$ws.Range("K5").Value = -0.015244513179135832
$ws.Range("B6").Value = 0.2667744730508082
$ws.Range("C6").Value = -0.03137900777207581
$ws.Range("D6").Value = 0.0
$ws.Range("E6").Value = -0.0023577010802306544
$ws.Range("F6").Value = 0.0006501898722127384
$ws.Range("G6").Value = 0.0
$ws.Range("H6").Value = 0.001490174199432903
$ws.Range("I6").Value = -0.011476021907462282
$ws.Range("J6").Value = 0.0
$ws.Range("K6").Value = -0.0002012619878772881
$ws.Range("B7").Value = 0.22103534590635748
$ws.Range("C7").Value = 0.0
$ws.Range("D7").Value = -0.0330201192580935
$ws.Range("E7").Value = -0.00172549812021947
$ws.Range("F7").Value = 0.008773886617471149
$ws.Range("G7").Value = 0.0024233210980385816
$ws.Range("H7").Value = 0.0
$ws.Range("I7").Value = 0.0
$ws.Range("J7").Value = 0.0
$ws.Range("K7").Value = 0.002264636453838287

# --- Column width adjustments (D, F, G, J) ---
$ws.Columns.Item(4).ColumnWidth = 14.833333333333332
$ws.Columns.Item(6).ColumnWidth = 15.0
$ws.Columns.Item(7).ColumnWidth = 14.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
